# Applies the diff: inserts two new data rows (623 & 624) into the
# "Fruta, Feria Lagunitas de Puerto Montt - Plátano" sheet, pushing the
# existing rows 623-692 down to 625-694.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new rows at 623, shifting existing rows 623:692 down to 625:694.
$ws.Rows("623:624").Insert()

# 2. Populate the two brand-new rows with their data (constant columns
#    A,B,C,E,F,G,H,I,J,Q,R,T are identical to every other row in this sheet).
$ws.Cells.Item(623,1).Value2  = 4
$ws.Cells.Item(623,2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(623,3).Value2  = "Los Lagos"
$ws.Cells.Item(623,4).Value2  = 44918
$ws.Cells.Item(623,5).Value2  = 10
$ws.Cells.Item(623,6).Value2  = "Fruta"
$ws.Cells.Item(623,7).Value2  = 100108
$ws.Cells.Item(623,8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(623,9).Value2  = 100108006
$ws.Cells.Item(623,10).Value2 = "Plátano"
$ws.Cells.Item(623,11).Value2 = "Sin especificar"
$ws.Cells.Item(623,12).Value2 = "Pintón"
$ws.Cells.Item(623,13).Value2 = 600
$ws.Cells.Item(623,14).Value2 = 25000
$ws.Cells.Item(623,15).Value2 = 25000
$ws.Cells.Item(623,16).Value2 = 25000
$ws.Cells.Item(623,17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(623,18).Value2 = "Ecuador"
$ws.Cells.Item(623,19).Value2 = 1250
$ws.Cells.Item(623,20).Value2 = 20

$ws.Cells.Item(624,1).Value2  = 4
$ws.Cells.Item(624,2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(624,3).Value2  = "Los Lagos"
$ws.Cells.Item(624,4).Value2  = 44918
$ws.Cells.Item(624,5).Value2  = 10
$ws.Cells.Item(624,6).Value2  = "Fruta"
$ws.Cells.Item(624,7).Value2  = 100108
$ws.Cells.Item(624,8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(624,9).Value2  = 100108006
$ws.Cells.Item(624,10).Value2 = "Plátano"
$ws.Cells.Item(624,11).Value2 = "Sin especificar"
$ws.Cells.Item(624,12).Value2 = "Primera Pintón"
$ws.Cells.Item(624,13).Value2 = 1200
$ws.Cells.Item(624,14).Value2 = 26000
$ws.Cells.Item(624,15).Value2 = 27000
$ws.Cells.Item(624,16).Value2 = 26500
$ws.Cells.Item(624,17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(624,18).Value2 = "Ecuador"
$ws.Cells.Item(624,19).Value2 = 1325
$ws.Cells.Item(624,20).Value2 = 20

# 3. The source diff leaves the "Precio máximo" (column O) of the row
#    that ends up at row 681 unchanged at its original (pre-shift) value
#    of 19000, instead of taking on the rest of that shifted row's new
#    25000 value, so restore it explicitly here.
$ws.Cells.Item(681,15).Value2 = 19000
